$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Fixed bom for Digikey"
#  - BAS16J135            -> BAS16J,135           (D1 D2 row)
#  - BKP2125HS221-T       -> MH2029-221Y          (FB1 FB2 row), with
#    left-to-right mark characters around the value as in the source data
#  - Remove the H1 / H2 / J1 placeholder rows (their "~" MPN is replaced by
#    folding that hardware into the FB1 FB2 / MH2029-221Y row above)
#  - MAX6008AEUR-T        -> MAX6008AEUR+TCT-ND   (U1 row)
# ---------------------------------------------------------------------------

# Fix the D1 D2 MPN (row 8)
$ws.Range("B8").Value = "BAS16J,135"

# Fix the FB1 FB2 MPN (row 11) - includes U+200E LEFT-TO-RIGHT MARK wrapping
$lrm = [char]0x200E
$ws.Range("B11").Value = "$lrm" + "MH2029-221Y" + "$lrm"

# Remove the H1 / H2 / J1 rows (rows 12, 13, 14)
$ws.Rows("12:14").Delete()

# Fix the U1 MPN (now row 32 after the row deletion above)
$lrm2 = [char]0x200E
$ws.Range("B32").Value = "MAX6008AEUR+TCT-ND" + "$lrm2"

# Cosmetic touch-ups from the same edit session
$ws.Range("B1").Font.Name = $ws.Range("A1").Font.Name()
$ws.Range("B1").Font.Size = $ws.Range("A1").Font.Size()
$ws.Range("B1").Font.Color = $ws.Range("A1").Font.Color()

$ws.Columns.Item(1).ColumnWidth = 75.9

$ws.Range("B11").Select()
